# Reassign the per-row "date / volume / price / origin" facts to different
# rows (a row permutation) while leaving all other columns untouched.
#
# Target row -> source row (the row whose D,M,N,O,P,R,S values should end up
# in the target row):
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 7
    3  = 10
    4  = 15
    5  = 6
    6  = 17
    7  = 3
    8  = 13
    9  = 5
    10 = 14
    11 = 11
    12 = 4
    13 = 8
    14 = 2
    15 = 9
    16 = 12
    17 = 16
}

$cols = @("D", "M", "N", "O", "P", "R", "S")

# Snapshot the original values of the columns we are going to move, keyed by
# row number, before we start overwriting anything.
$original = @{}
for ($r = 2; $r -le 17; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowVals
}

# Now write back the permuted values.
for ($r = 2; $r -le 17; $r++) {
    $srcRow = $mapping[$r]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $srcVals[$col]
    }
}
